$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Range("F4").Value = "8.0.52"
$ws.Range("F4").Select()
